# Update and Delete test passed for Zones.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")

# Rows 2-10 correspond to Zone test results (Z01-Z09).
# Mark "Update Test Passed" (column D) and "Delete Test Passed" (column E) as TRUE.
$ws.Range("D2:E10").Value = $true

# Scroll the sheet view back to the top (no longer showing topLeftCell A11).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

$ws.Range("C21").Select()
